# Apply the content edit described by the commit "Cập nhật file tài liệu 2".
#
# The only true content-level change in the diff is the rewording of the
# bullet about personal leave: the old one-sentence bullet is replaced by a
# longer sentence about notifying the management department at least 24h in
# advance. (The remaining hunks in the diff only move
# <w:lastRenderedPageBreak/> markers around, which are just Word's cached
# pagination bookkeeping following the reflow caused by this text edit --
# they carry no visible/textual content and are regenerated by Word itself
# whenever the document is repaginated, so there is nothing else to "edit"
# there.)

$d = $word.ActiveDocument

$old = "Đối với nghỉ việc cá nhân, cần sự chấp thuận từ quản lý."
$new = "Đối với cá nhân, nếu cần đổi ca làm việc hoặc xin nghỉ phép, nhân viên phải thông báo cho bộ phận quản lý ít nhất trước 24h."

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the target sentence to replace."
}

Write-Output "Replaced sentence: $found"
